$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hours update: two new timesheet rows for 1/26/2010 (Group Meeting + Weekly Meeting)
# Format the date cells as text first so "1/26/2010" is stored as a literal
# string (matching the rest of column A) instead of being auto-converted to
# a date serial number, then drop the temporary number format again so the
# cells keep the sheet's default (unstyled) look.
$ws.Range("A26:A27").NumberFormat = "@"

$ws.Cells.Item(26, 1).Value = "1/26/2010"
$ws.Cells.Item(26, 2).Value = 2
$ws.Cells.Item(26, 3).Value = "Group Meeting"

$ws.Cells.Item(27, 1).Value = "1/26/2010"
$ws.Cells.Item(27, 2).Value = 1
$ws.Cells.Item(27, 3).Value = "Weekly Meeting"

$ws.Range("A26:A27").ClearFormats()

# Update the selection to match the post-edit state
$ws.Range("C28").Select()
